# issue #5: add legislator_id, name, date into dataframe
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 56

# --- Header row (row 1): new columns H/I/J -----------------------------
$ws.Cells.Item(1, 8).Value  = "date"
$ws.Cells.Item(1, 9).Value  = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Match the look of the existing header cells (bold, centered, bordered)
$headerSample = $ws.Cells.Item(1, 7)
foreach ($col in 8, 9, 10) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Font.Bold = $headerSample.Font.Bold
    $headerCell.HorizontalAlignment = $headerSample.HorizontalAlignment
    $headerCell.VerticalAlignment = $headerSample.VerticalAlignment
    $headerCell.Borders.LineStyle = $headerSample.Borders.LineStyle
}

# --- Data rows (2..56): same date / legislator name / legislator id ----
$date = "2011-11-22"
$legislatorName = "馬文君"
$legislatorId = 1724

# Format column H as text first so the date-looking string is stored as
# literal text instead of being auto-converted into a date serial number.
$ws.Range("H2:H$lastRow").NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value  = $date
    $ws.Cells.Item($r, 9).Value  = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
